$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.217.03'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.010.66'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.67'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.645'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.34'
$ws.Range('E7').Value = '  +17.76%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '59.22'
$ws.Range('E9').Value = '  -2.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.370'
$ws.Range('E10').Value = '  +3.21%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0749'
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.948'
$ws.Range('E13').Value = '  +2.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.98'
$ws.Range('E14').Value = '  +3.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.293.91'
$ws.Range('E15').Value = '  -1.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.45'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.64'
$ws.Range('E17').Value = '  +17.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.989.33'
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '36.121.96'
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.27'
$ws.Range('E20').Value = '  +1.87%  '
$ws.Range('E21').Value = '  +1.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.27'
$ws.Range('E22').Value = '  +3.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.52'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.64'
$ws.Range('E25').Value = '  +18.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('E26').Value = '  -3.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.74'
$ws.Range('E27').Value = '  +7.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '166.14'
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.70'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.13'
$ws.Range('E31').Value = '  +4.49%  '
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('E33').Value = '  +15.61%  '
$ws.Range('E34').Value = '  +3.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.50'
$ws.Range('E35').Value = '  +3.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.48'
$ws.Range('E36').Value = '  +13.74%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('E39').Value = '  +18.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.25'
$ws.Range('E40').Value = '  +3.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0966'
$ws.Range('E41').Value = '  +9.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.90'
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  +2.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.99'
$ws.Range('E44').Value = '  +9.52%  '
$ws.Range('E45').Value = '  +3.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '94.64'
$ws.Range('E46').Value = '  +2.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.85'
$ws.Range('E47').Value = '  +7.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.373.75'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.93'
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.35'
$ws.Range('E50').Value = '  +5.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '47.26'
$ws.Range('E51').Value = '  +6.41%  '
